# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect freshly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 19
$ws1.Range("F3").Value = 12252
$ws1.Range("F8").Value = 31
$ws1.Range("F9").Value = 2623
$ws1.Range("F10").Value = 1138
$ws1.Range("F11").Value = 217
$ws1.Range("F13").Value = 5373
$ws1.Range("F15").Value = 220
$ws1.Range("F16").Value = 570
$ws1.Range("F17").Value = 11522
$ws1.Range("F18").Value = 11643

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 19
$ws4.Range("F3").Value = 12252
$ws4.Range("F8").Value = 31
$ws4.Range("F9").Value = 2623
$ws4.Range("F11").Value = 1138
$ws4.Range("F12").Value = 217
$ws4.Range("F14").Value = 5373
$ws4.Range("F16").Value = 220
$ws4.Range("F17").Value = 570
$ws4.Range("F18").Value = 11522
$ws4.Range("F19").Value = 11643
